# Apply cell updates per the target diff (scrape snapshot refresh: view counts / prices / statuses,
# plus one fully-replaced event row in the "全部类型" aggregate sheet).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('展览')
$ws.Range('F2').Value = 221
$ws.Range('F3').Value = 54884
$ws.Range('F5').Value = 389
$ws.Range('F6').Value = 350
$ws.Range('G6').Value = '不可售'
$ws.Range('F7').Value = 887
$ws.Range('F8').Value = 778
$ws.Range('F9').Value = 416
$ws.Range('F10').Value = 3095
$ws.Range('F11').Value = 921
$ws.Range('G13').Value = '不可售'
$ws.Range('F14').Value = 1114
$ws.Range('F18').Value = 428
$ws.Range('F19').Value = 1324
$ws.Range('F20').Value = 109
$ws.Range('F21').Value = 41
$ws.Range('F22').Value = 199
$ws.Range('F23').Value = 383
$ws.Range('F24').Value = 42
$ws.Range('F27').Value = 71
$ws.Range('F28').Value = 63
$ws.Range('F29').Value = 5222
$ws.Range('F30').Value = 39
$ws.Range('F31').Value = 5151
$ws.Range('F32').Value = 9148
$ws.Range('F34').Value = 156
$ws.Range('F37').Value = 441
$ws.Range('F38').Value = 126
$ws.Range('F40').Value = 4240
$ws.Range('F41').Value = 263
$ws = $wb.Worksheets.Item('演出')
$ws.Range('F4').Value = 103
$ws = $wb.Worksheets.Item('本地生活')
$ws.Range('F3').Value = 578
$ws.Range('F4').Value = 139
$ws.Range('F5').Value = 48
$ws = $wb.Worksheets.Item('全部类型')
$ws.Range('F3').Value = 578
$ws.Range('F4').Value = 221
$ws.Range('B5').NumberFormat = "@"
$ws.Range('B5').Value = '2024-10-04'
$ws.Range('B5').ClearFormats()
$ws.Range('C5').Value = '杭州·创世次元第五人格同人only展'
$ws.Range('D5').Value = '小河路与桥弄街交叉口东北50米 桥西历史文化街区'
$ws.Range('E5').Value = '2024.10.04 10:00-10.05 17:00'
$ws.Range('F5').Value = 1347
$ws.Range('G5').Value = 85
$ws.Range('H5').Value = 'https://show.bilibili.com/platform/detail.html?id=92141'
$ws.Range('I5').Value = '//i1.hdslb.com/bfs/openplatform/202409/MMF3dkAw1725550270634.jpeg'
$ws.Range('F6').Value = 887
$ws.Range('F7').Value = 778
$ws.Range('F8').Value = 416
$ws.Range('F9').Value = 3095
$ws.Range('F10').Value = 921
$ws.Range('F11').Value = 103
$ws.Range('F12').Value = 48
$ws.Range('F14').Value = 1114
$ws.Range('F18').Value = 428
$ws.Range('F20').Value = 1324
$ws.Range('F22').Value = 109
$ws.Range('F23').Value = 199
$ws.Range('F25').Value = 383
$ws.Range('F27').Value = 63
$ws.Range('F28').Value = 5222
$ws.Range('F29').Value = 5151
$ws.Range('F30').Value = 9148
$ws.Range('F33').Value = 156
$ws.Range('F35').Value = 441
$ws.Range('F39').Value = 126
$ws.Range('F41').Value = 4240
$ws.Range('F47').Value = 263
